$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.729.93'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").Value = '2.310.89'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.10'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.979'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").Value = '2.660.41'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").Value = '2.308.06'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").Value = '42.656.47'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +35.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0887'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.72%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0354'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("E40").Value = '  -4.99%  '
$ws.Range("E41").Value = '  +12.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.69%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '116.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '82.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.56%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.638.52'
$ws.Range("E49").Value = '  +4.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.97%  '
